$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-49
$ws.Range("D2").Value = "'28.311.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.34%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.578.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.17%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +3.11%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'211.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.10%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -1.17%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +3.76%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'46.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.30%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'23.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.37%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -1.22%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0595"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.84%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.22%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.805.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.07%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.574.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.09%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -1.63%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'28.360.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.86%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'62.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.97%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'228.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.62%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.62%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.62%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +2.07%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -4.61%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.86%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +4.09%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'151.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.18%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'15.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.58%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.67%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.63%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +2.70%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.16%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.68%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.87%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.62%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.389.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.15%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.77%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.46%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +2.92%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.84%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.541"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.34%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.802"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.98%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +2.33%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.76%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.71%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +1.78%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'62.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.25%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.715.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.31%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'85.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.41%  "
$ws.Range("E49").Style = "Normal"

# Row 50/51: Cronos and BabyDogeCoin swap order, with new Price/Volume values
$ws.Range("B50").Value = "'Cronos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0520"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.87%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'BabyDogeCoin"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0₆01000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.98%  "
$ws.Range("E51").Style = "Normal"
